# Updated cryptos list -- apply price/volume changes row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain a text cell even when the string looks numeric
    # (e.g. "21.35", "0.0780"), matching the source data which stores every
    # Price/Volume cell as literal text -- and restore the original (default)
    # number format/style afterwards so no stray style gets attached.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.092.04"
$ws.Range("E2").Value = "  -0.72%  "

Set-TextValue $ws.Range("D3") "2.604.46"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "589.93"
$ws.Range("E5").Value = "  -1.81%  "

Set-TextValue $ws.Range("D6") "150.08"
$ws.Range("E6").Value = "  -2.49%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -0.57%  "

Set-TextValue $ws.Range("D9") "2.603.62"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("E13").Value = "  -3.18%  "

Set-TextValue $ws.Range("D14") "27.14"
$ws.Range("E14").Value = "  -2.84%  "

Set-TextValue $ws.Range("D15") "3.076.93"
$ws.Range("E15").Value = "  -0.49%  "

Set-TextValue $ws.Range("D16") "0.0000181"
$ws.Range("E16").Value = "  -2.97%  "

Set-TextValue $ws.Range("D17") "66.921.53"
$ws.Range("E17").Value = "  -1.00%  "

Set-TextValue $ws.Range("D18") "2.606.31"
$ws.Range("E18").Value = "  -0.43%  "

Set-TextValue $ws.Range("D19") "367.50"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("E20").Value = "  -2.05%  "

$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("E22").Value = "  -0.54%  "

Set-TextValue $ws.Range("D23") "4.74"
$ws.Range("E23").Value = "  -4.72%  "

$ws.Range("E24").Value = "  -3.39%  "

Set-TextValue $ws.Range("D25") "73.51"
$ws.Range("E25").Value = "  +4.98%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -0.79%  "

Set-TextValue $ws.Range("D28") "2.737.27"
$ws.Range("E28").Value = "  -0.37%  "

Set-TextValue $ws.Range("D29") "583.60"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -0.12%  "

Set-TextValue $ws.Range("D31") "0.0₃0986"
$ws.Range("E31").Value = "  -6.70%  "

$ws.Range("E32").Value = "  -5.16%  "

Set-TextValue $ws.Range("D33") "7.65"
$ws.Range("E33").Value = "  -3.81%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("E36").Value = "  -4.31%  "

$ws.Range("E37").Value = "  -2.95%  "

Set-TextValue $ws.Range("D38") "156.55"
$ws.Range("E38").Value = "  +0.31%  "

Set-TextValue $ws.Range("D39") "19.02"
$ws.Range("E39").Value = "  -2.09%  "

Set-TextValue $ws.Range("D40") "0.365"
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("E42").Value = "  -3.31%  "

Set-TextValue $ws.Range("D43") "2.56"
$ws.Range("E43").Value = "  -4.44%  "

$ws.Range("E44").Value = "  +4.15%  "

Set-TextValue $ws.Range("D45") "0.999"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("E46").Value = "  -2.64%  "

Set-TextValue $ws.Range("D47") "0.0₆0286"
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("E48").Value = "  -1.44%  "

Set-TextValue $ws.Range("D51") "21.35"
$ws.Range("E51").Value = "  +1.64%  "

# Rows 49 and 50 swap places in the ranking (Optimism <-> Cronos) and their
# price/volume figures are refreshed as well.
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0780"
$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextValue $ws.Range("D50") "1.68"
$ws.Range("E50").Value = "  -3.69%  "

